$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.930.84"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "2.263.39"
$ws.Range("E3").Value = "  -0.43%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.65"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.52%  "
$ws.Range("E7").Value = "  -0.77%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.489"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.89"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.67%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0788"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.29%  "
$ws.Range("E12").Value = "  -0.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.63"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.91%  "
$ws.Range("D14").Value = "2.616.52"
$ws.Range("E14").Value = "  -0.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.36"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.17%  "
$ws.Range("D16").Value = "2.262.11"
$ws.Range("E16").Value = "  -0.62%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.791"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.00%  "
$ws.Range("D18").Value = "41.845.03"
$ws.Range("E18").Value = "  +0.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.30"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.97%  "
$ws.Range("D20").Value = "0.0₃0900"
$ws.Range("E20").Value = "  -1.81%  "
$ws.Range("E21").Value = "  -0.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.06"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.80%  "
$ws.Range("E24").Value = "  -0.89%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("B26").Value = "ImmutableX"
$ws.Range("C26").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.92"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.85%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.63"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.63%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "36.58"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.13%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.12"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.99%  "
$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.49"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.93%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "160.12"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.20"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.25%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.14"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0735"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.12"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.88%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.39"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.65%  "
$ws.Range("E38").Value = "  -0.89%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.82"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.50%  "
$ws.Range("E40").Value = "  -2.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.70%  "
$ws.Range("E42").Value = "  +5.83%  "
$ws.Range("D43").Value = "1.975.00"
$ws.Range("E43").Value = "  -1.46%  "
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.85"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.83%  "
$ws.Range("E46").Value = "  +0.64%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.85"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.08%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.07"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.74%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "72.95"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.37%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.50"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "90.72"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.31%  "
